$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-24 21:18:18"
$ws.Range("O2").Value = "5.8 °C"
$ws.Range("E3").Value = "2026-02-24 21:18:20"
$ws.Range("E4").Value = "2026-02-24 21:18:23"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "72%"
$ws.Range("J4").Value = "1019.7 hPa"
$ws.Range("O4").Value = "12.9 °C"
$ws.Range("E5").Value = "2026-02-24 21:18:25"
$ws.Range("O5").Value = "5.9 °C"
$ws.Range("E6").Value = "2026-02-24 21:18:27"
$ws.Range("E7").Value = "2026-02-24 21:18:29"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "72%"
$ws.Range("O7").Value = "14.1 °C"
$ws.Range("E8").Value = "2026-02-24 21:18:32"
$ws.Range("E9").Value = "2026-02-24 21:18:35"
$ws.Range("O9").Value = "11.8 °C"
$ws.Range("E10").Value = "2026-02-24 21:18:37"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "77%"
$ws.Range("O10").Value = "11.2 °C"
$ws.Range("E11").Value = "2026-02-24 21:18:39"
$ws.Range("E12").Value = "2026-02-24 21:18:42"
$ws.Range("O12").Value = "10.6 °C"
$ws.Range("E13").Value = "2026-02-24 21:18:44"
$ws.Range("J13").Value = "1023.2 hPa"
$ws.Range("E14").Value = "2026-02-24 21:18:46"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "85%"
$ws.Range("N14").Value = "6.4 °C 20:54 TU"
$ws.Range("O14").Value = "11.6 °C"
$ws.Range("E15").Value = "2026-02-24 21:18:48"
$ws.Range("O15").Value = "11.9 °C"
$ws.Range("E16").Value = "2026-02-24 21:18:50"
$ws.Range("E17").Value = "2026-02-24 21:18:51"
$ws.Range("E18").Value = "2026-02-24 21:18:52"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "76%"
$ws.Range("E19").Value = "2026-02-24 21:18:53"
$ws.Range("O19").Value = "12.6 °C"
$ws.Range("E20").Value = "2026-02-24 21:18:54"
$ws.Range("E21").Value = "2026-02-24 21:18:55"
$ws.Range("E22").Value = "2026-02-24 21:18:58"
$ws.Range("O22").Value = "3.5 °C"
$ws.Range("E23").Value = "2026-02-24 21:19:00"
$ws.Range("K23").Value = "16.1 MJ/m2"
$ws.Range("E24").Value = "2026-02-24 21:19:03"
$ws.Range("J24").Value = "1021.1 hPa"
$ws.Range("L24").Value = "10.8 km/h - 75º 20:51 TU"
$ws.Range("O24").Value = "9.8 °C"
$ws.Range("E25").Value = "2026-02-24 21:19:05"
$ws.Range("O25").Value = "6.7 °C"
$ws.Range("E26").Value = "2026-02-24 21:19:07"
$ws.Range("O26").Value = "11.7 °C"
$ws.Range("E27").Value = "2026-02-24 21:19:10"
$ws.Range("O27").Value = "6.3 °C"
$ws.Range("E28").Value = "2026-02-24 21:19:12"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "71%"
$ws.Range("O28").Value = "11.6 °C"
$ws.Range("E29").Value = "2026-02-24 21:19:14"
$ws.Range("E30").Value = "2026-02-24 21:19:17"
$ws.Range("O30").Value = "13.1 °C"
$ws.Range("E31").Value = "2026-02-24 21:19:19"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "60%"
$ws.Range("E32").Value = "2026-02-24 21:19:22"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "68%"
$ws.Range("O32").Value = "7.2 °C"
$ws.Range("E33").Value = "2026-02-24 21:19:24"
$ws.Range("O33").Value = "8.6 °C"
$ws.Range("E34").Value = "2026-02-24 21:19:27"
$ws.Range("O34").Value = "4.7 °C"
$ws.Range("E35").Value = "2026-02-24 21:19:29"
$ws.Range("E36").Value = "2026-02-24 21:19:32"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "81%"
$ws.Range("J36").Value = "1019.8 hPa"
$ws.Range("E37").Value = "2026-02-24 21:19:34"
$ws.Range("O37").Value = "8.7 °C"
$ws.Range("E38").Value = "2026-02-24 21:19:36"
$ws.Range("O38").Value = "12.0 °C"
$ws.Range("E39").Value = "2026-02-24 21:19:39"
$ws.Range("E40").Value = "2026-02-24 21:19:41"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "66%"
$ws.Range("O40").Value = "8.6 °C"
$ws.Range("E41").Value = "2026-02-24 21:19:44"
$ws.Range("E42").Value = "2026-02-24 21:19:46"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "87%"
$ws.Range("E43").Value = "2026-02-24 21:19:48"
$ws.Range("E44").Value = "2026-02-24 21:19:51"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "40%"
$ws.Range("O44").Value = "2.8 °C"
$ws.Range("E45").Value = "2026-02-24 21:19:53"
$ws.Range("O45").Value = "10.1 °C"
$ws.Range("E46").Value = "2026-02-24 21:19:56"
$ws.Range("J46").Value = "1021.1 hPa"
